$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = 0.04908414824669194
    "D2" = 0.02285065723522273
    "E2" = 0.08923344112572096
    "F2" = 0.9892913280850166
    "G2" = 0.8521200284607033
    "H2" = 0.8602198199418893
    "K2" = 1.215002793583153
    "M2" = 0.3931538099406495
    "C3" = 0.04357890026741984
    "D3" = 0.02310382918422427
    "E3" = 0.0846545358746873
    "F3" = 0.9733567836959907
    "G3" = 0.8357737081999375
    "H3" = 0.8593017098962719
    "K3" = 1.07158942560136
    "M3" = 0.3531198535792441
    "C4" = 0.04021604835791948
    "D4" = 0.02326941884447642
    "E4" = 0.08192774140342607
    "F4" = 0.9644558733761244
    "G4" = 0.8265939783176179
    "H4" = 0.8593748282463167
    "K4" = 0.9835886444212463
    "M4" = 0.3286707486100084
    "C5" = 0.03884994313415291
    "D5" = 0.02333944588709258
    "E5" = 0.08083753496028123
    "F5" = 0.9610492454698516
    "G5" = 0.8230668815005373
    "H5" = 0.8595640725918088
    "K5" = 0.9477417793427207
    "M5" = 0.3187402583377832
    "C6" = 0.03862335857068899
    "D6" = 0.02335122766109876
    "E6" = 0.08065776547421777
    "F6" = 0.9604968563629228
    "G6" = 0.8224940658713393
    "H6" = 0.8596051065646435
    "K6" = 0.9417902969294971
    "M6" = 0.3170932712149153
    "C7" = 0.04019760735455691
    "D7" = 0.0232703529399263
    "E7" = 0.08191295393667986
    "F7" = 0.9644090392851297
    "G7" = 0.826545547678279
    "H7" = 0.8593767357882598
    "K7" = 0.98310514341199
    "M7" = 0.3285366906524843
    "C8" = 0.04718227329689739
    "D8" = 0.02293584467701315
    "E8" = 0.08763687426301203
    "F8" = 0.9836129490058454
    "G8" = 0.8463048750141837
    "H8" = 0.8597706599350232
    "K8" = 1.165542376577605
    "M8" = 0.3793224379474367
    "C9" = 0.06102192432767595
    "D9" = 0.02236046429930205
    "E9" = 0.09954780427774779
    "F9" = 1.028346988361932
    "G9" = 0.8919357459078583
    "H9" = 0.8656286878378694
    "K9" = 1.523754787574262
    "M9" = 0.479985702874032
    "C10" = 0.07128452510127659
    "D10" = 0.02198707500239649
    "E10" = 0.1087382133640347
    "F10" = 1.065627061577047
    "G10" = 0.9297762795294204
    "H10" = 0.873079768908525
    "K10" = 1.787256869934481
    "M10" = 0.5546409417873122
    "C11" = 0.07597550071029957
    "D11" = 0.02182798557932841
    "E11" = 0.1130192139581538
    "F11" = 1.083567249107773
    "G11" = 0.9479540824455057
    "H11" = 0.8771631241645252
    "K11" = 1.907213667732037
    "M11" = 0.5887652705097537
    "C12" = 0.07775519957745303
    "D12" = 0.0217692973020398
    "E12" = 0.1146550915222946
    "F12" = 1.090503543150191
    "G12" = 0.9549781940307014
    "H12" = 0.8788099631826469
    "K12" = 1.952651385875299
    "M12" = 0.6017114241057158
    "C13" = 0.07737175999615431
    "D13" = 0.02178186757305767
    "E13" = 0.1143021142898135
    "F13" = 1.089003313993445
    "G13" = 0.9534591426793497
    "H13" = 0.8784508018107715
    "K13" = 1.942864996139406
    "M13" = 0.5989221617989102
    "C14" = 0.07612185034211905
    "D14" = 0.02182312604346848
    "E14" = 0.1131535008986972
    "F14" = 1.084135032242628
    "G14" = 0.9485291318165139
    "H14" = 0.8772965906634056
    "K14" = 1.91095160195448
    "M14" = 0.5898298740521284
    "C15" = 0.07535668093390768
    "D15" = 0.02184860083440654
    "E15" = 0.1124518738309135
    "F15" = 1.081171705719072
    "G15" = 0.9455277226867622
    "H15" = 0.876602722566048
    "K15" = 1.891405378194463
    "M15" = 0.5842637281969303
    "C16" = 0.07097841921404324
    "D16" = 0.02199768912540279
    "E16" = 0.1084604860625191
    "F16" = 1.064474513035478
    "G16" = 0.9286078816316206
    "H16" = 0.8728269409355107
    "K16" = 1.779419186098949
    "M16" = 0.5524141674103618
    "C17" = 0.06829832079002074
    "D17" = 0.02209191253855991
    "E17" = 0.1060378141531615
    "F17" = 1.054483813071428
    "G17" = 0.9184763882624054
    "H17" = 0.8706888861224797
    "K17" = 1.710741885930588
    "M17" = 0.5329176879040176
    "C18" = 0.06675890786708294
    "D18" = 0.02214712034583854
    "E18" = 0.1046537715565137
    "F18" = 1.048829693184672
    "G18" = 0.9127396356836357
    "H18" = 0.8695243917223934
    "K18" = 1.671248836425889
    "M18" = 0.5217191651951509
    "C19" = 0.06623804834045188
    "D19" = 0.02216598653660284
    "E19" = 0.1041867650594952
    "F19" = 1.046931102809694
    "G19" = 0.9108127734731397
    "H19" = 0.8691412980757889
    "K19" = 1.657878598544983
    "M19" = 0.5179301559788598
    "C20" = 0.06858340288471254
    "D20" = 0.02208177740656225
    "E20" = 0.1062947349128223
    "F20" = 1.055537779184874
    "G20" = 0.9195455106517443
    "H20" = 0.8709097260511953
    "K20" = 1.71805184130983
    "M20" = 0.5349915303381181
    "C21" = 0.07648888799640474
    "D21" = 0.02181096515646885
    "E21" = 0.1134904732953572
    "F21" = 1.085561077833788
    "G21" = 0.9499733647831476
    "H21" = 0.8776328749330276
    "K21" = 1.920324995193937
    "M21" = 0.592499843414501
    "C22" = 0.08167502497930457
    "D22" = 0.02164304379029502
    "E22" = 0.1182794717553151
    "F22" = 1.106015765023443
    "G22" = 0.9706800033857235
    "H22" = 0.8826133585338027
    "K22" = 2.052596138808667
    "M22" = 0.6302250609795692
    "C23" = 0.07890527428631344
    "D23" = 0.0217318341230559
    "E23" = 0.115715495128363
    "F23" = 1.095021979650497
    "G23" = 0.9595527762336076
    "H23" = 0.8799012517788753
    "K23" = 1.981993770129009
    "M23" = 0.6100773940174236
    "C24" = 0.06845451288273807
    "D24" = 0.02208635626830002
    "E24" = 0.1061785537050852
    "F24" = 1.055061001924969
    "G24" = 0.9190618865276008
    "H24" = 0.8708096828345617
    "K24" = 1.714747040253087
    "M24" = 0.5340539142521123
    "C25" = 0.05726187007421402
    "D25" = 0.02250747300760203
    "E25" = 0.09624994052632729
    "F25" = 1.01547740891354
    "G25" = 0.8788421978759402
    "H25" = 0.8634946619776827
    "K25" = 1.426795614467665
    "M25" = 0.4526342922638378
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
